# Refresh crypto price/volume snapshot values (scraped data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.712.94"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.645.11"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D5").Value = "'213.45"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'23.19"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "1.878.20"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "1.644.40"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "'0.563"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "'64.29"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "27.690.01"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "'231.63"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'7.67"
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "'10.11"
$ws.Range("E23").Value = "  +7.24%  "
$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  -3.53%  "
$ws.Range("D25").Value = "'150.02"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'15.68"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'0.0487"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "1.446.83"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'3.16"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.884"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "'0.900"
$ws.Range("E40").Value = "  +14.16%  "
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'66.12"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").Value = "'2.25"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "1.787.26"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("D49").Value = "'86.52"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("E51").Value = "  -1.85%  "
